# Updates the cryptos price/volume table to the latest scrape.
# For D-column numeric-looking strings we prefix with a literal
# apostrophe (Excel's quote-prefix) so the cell stays text instead
# of being coerced to a Number (which would silently drop things
# like trailing zeros, e.g. "1.00" -> 1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '27.439.18'
$ws.Range('E2').Value = '  -0.60%  '

# Row 3
$ws.Range('D3').Value = '1.643.84'
$ws.Range('E3').Value = '  -1.32%  '

# Row 4
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  -0.01%  '

# Row 5
$ws.Range('D5').Value = '''212.19'
$ws.Range('E5').Value = '  -1.49%  '

# Row 6
$ws.Range('D6').Value = '''0.538'
$ws.Range('E6').Value = '  +4.81%  '

# Row 7
$ws.Range('E7').Value = '  -0.04%  '

# Row 8
$ws.Range('D8').Value = '''23.14'
$ws.Range('E8').Value = '  -1.99%  '

# Row 9
$ws.Range('E9').Value = '  -2.24%  '

# Row 10
$ws.Range('D10').Value = '''0.0609'
$ws.Range('E10').Value = '  -1.97%  '

# Row 11
$ws.Range('E11').Value = '  +0.72%  '

# Row 12
$ws.Range('D12').Value = '1.876.67'
$ws.Range('E12').Value = '  -1.26%  '

# Row 13
$ws.Range('D13').Value = '1.659.30'
$ws.Range('E13').Value = '  -0.69%  '

# Row 14
$ws.Range('D14').Value = '''4.03'
$ws.Range('E14').Value = '  -2.95%  '

# Row 15
$ws.Range('D15').Value = '''0.557'
$ws.Range('E15').Value = '  +0.10%  '

# Row 16
$ws.Range('D16').Value = '''64.34'
$ws.Range('E16').Value = '  -2.84%  '

# Row 17
$ws.Range('D17').Value = '27.404.30'

# Row 18
$ws.Range('D18').Value = '''228.52'
$ws.Range('E18').Value = '  -7.76%  '

# Row 19
$ws.Range('D19').Value = '0.0₃0720'
$ws.Range('E19').Value = '  -1.50%  '

# Row 20
$ws.Range('D20').Value = '''7.50'
$ws.Range('E20').Value = '  -0.61%  '

# Row 22
$ws.Range('D22').Value = '''4.33'
$ws.Range('E22').Value = '  -3.65%  '

# Row 23
$ws.Range('D23').Value = '''9.32'
$ws.Range('E23').Value = '  +0.30%  '

# Row 24
$ws.Range('D24').Value = '''2.03'
$ws.Range('E24').Value = '  -0.47%  '

# Row 25
$ws.Range('D25').Value = '''148.11'
$ws.Range('E25').Value = '  +1.41%  '

# Row 26
$ws.Range('D26').Value = '''0.114'
$ws.Range('E26').Value = '  +2.51%  '

# Row 27
$ws.Range('E27').Value = '  -3.08%  '

# Row 28
$ws.Range('E28').Value = '  -0.02%  '

# Row 29
$ws.Range('D29').Value = '''15.54'
$ws.Range('E29').Value = '  -5.31%  '

# Row 30
$ws.Range('E30').Value = '  -4.86%  '

# Row 31
$ws.Range('E31').Value = '  -3.79%  '

# Row 32
$ws.Range('D32').Value = '''3.27'
$ws.Range('E32').Value = '  -2.25%  '

# Row 33
$ws.Range('D33').Value = '''3.11'
$ws.Range('E33').Value = '  +0.00%  '

# Row 34
$ws.Range('D34').Value = '1.413.61'
$ws.Range('E34').Value = '  -4.33%  '

# Row 35
$ws.Range('E35').Value = '  -0.10%  '

# Row 36
$ws.Range('E36').Value = '  -0.23%  '

# Row 37
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '''0.563'
$ws.Range('E37').Value = '  -1.78%  '

# Row 38
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').Value = '''0.880'
$ws.Range('E38').Value = '  -6.13%  '

# Row 39
$ws.Range('E39').Value = '  -3.22%  '

# Row 40
$ws.Range('E40').Value = '  +0.77%  '

# Row 41
$ws.Range('E41').Value = '  -0.03%  '

# Row 42
$ws.Range('E42').Value = '  -1.55%  '

# Row 43
$ws.Range('D43').Value = '''5.49'
$ws.Range('E43').Value = '  +1.42%  '

# Row 44
$ws.Range('E44').Value = '  +0.16%  '

# Row 45
$ws.Range('D45').Value = '''64.65'
$ws.Range('E45').Value = '  -6.89%  '

# Row 46
$ws.Range('D46').Value = '''0.790'
$ws.Range('E46').Value = '  +0.33%  '

# Row 47
$ws.Range('D47').Value = '1.787.36'
$ws.Range('E47').Value = '  -1.13%  '

# Row 48
$ws.Range('D48').Value = '''1.64'
$ws.Range('E48').Value = '  -3.54%  '

# Row 49
$ws.Range('D49').Value = '''87.52'
$ws.Range('E49').Value = '  -1.95%  '

# Row 50
$ws.Range('E50').Value = '  -3.45%  '

# Row 51
$ws.Range('E51').Value = '  -3.46%  '
